$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Update D2 text value (shared string "autotest10" -> "autotest12")
$ws.Range("D2").Value = "autotest12"

# Update B2 numeric value
$ws.Range("B2").Value = 48349402

# Change the active/selected cell on the sheet from B2 to D2
$ws.Activate()
$ws.Range("D2").Select()
